$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue 'D2' '28.407.96'
Set-TextValue 'E2' '  -0.31%  '
Set-TextValue 'D3' '1.863.82'
Set-TextValue 'E3' '  +0.14%  '
Set-TextValue 'E4' '  +0.10%  '
Set-TextValue 'D5' '324.81'
Set-TextValue 'E5' '  -0.47%  '
Set-TextValue 'D6' '1.008'
Set-TextValue 'E6' '  +0.20%  '
Set-TextValue 'D7' '0.4562'
Set-TextValue 'E7' '  -1.78%  '
Set-TextValue 'D8' '0.3837'
Set-TextValue 'E8' '  -1.40%  '
Set-TextValue 'D9' '0.07816'
Set-TextValue 'E9' '  -0.89%  '
Set-TextValue 'D10' '0.9861'
Set-TextValue 'E10' '  +1.37%  '
Set-TextValue 'D11' '21.55'
Set-TextValue 'E11' '  -3.15%  '
Set-TextValue 'D12' '1.829.52'
Set-TextValue 'E12' '  +1.06%  '
Set-TextValue 'B13' 'Polkadot'
Set-TextValue 'C13' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D13' '5.637'
Set-TextValue 'E13' '  -1.42%  '
Set-TextValue 'B14' 'Chainlink'
Set-TextValue 'C14' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D14' '6.894'
Set-TextValue 'E14' '  -0.56%  '
Set-TextValue 'D15' '0.06924'
Set-TextValue 'E15' '  +0.21%  '
Set-TextValue 'D16' '86.77'
Set-TextValue 'E16' '  -2.40%  '
Set-TextValue 'D17' '1.009'
Set-TextValue 'E17' '  +0.15%  '
Set-TextValue 'D18' '0.000009949'
Set-TextValue 'E18' '  -0.48%  '
Set-TextValue 'D19' '16.67'
Set-TextValue 'E19' '  -1.02%  '
Set-TextValue 'D20' '1.006'
Set-TextValue 'E20' '  +0.23%  '
Set-TextValue 'D21' '28.423.89'
Set-TextValue 'E21' '  -0.26%  '
Set-TextValue 'D22' '5.247'
Set-TextValue 'E22' '  -1.51%  '
Set-TextValue 'E23' '  -1.61%  '
Set-TextValue 'D24' '2.097'
Set-TextValue 'E24' '  -0.51%  '
Set-TextValue 'D25' '2.084.00'
Set-TextValue 'E25' '  -0.15%  '
Set-TextValue 'D26' '153.78'
Set-TextValue 'D27' '19.09'
Set-TextValue 'E27' '  -0.97%  '
Set-TextValue 'D28' '5.633'
Set-TextValue 'E28' '  -2.44%  '
Set-TextValue 'D29' '117.56'
Set-TextValue 'E29' '  -1.35%  '
Set-TextValue 'D30' '1.914'
Set-TextValue 'E30' '  -3.77%  '
Set-TextValue 'D31' '0.09284'
Set-TextValue 'E31' '  -0.21%  '
Set-TextValue 'D32' '0.9075'
Set-TextValue 'E32' '  -2.97%  '
Set-TextValue 'D33' '5.252'
Set-TextValue 'E33' '  -1.14%  '
Set-TextValue 'E34' '  -1.51%  '
Set-TextValue 'E35' '  -1.37%  '
Set-TextValue 'D36' '0.05686'
Set-TextValue 'E36' '  -2.62%  '
Set-TextValue 'D37' '1.158'
Set-TextValue 'E37' '  +0.55%  '
Set-TextValue 'D38' '0.02051'
Set-TextValue 'E38' '  -3.24%  '
Set-TextValue 'D39' '7.663'
Set-TextValue 'E39' '  -1.88%  '
Set-TextValue 'D40' '0.5554'
Set-TextValue 'E40' '  -1.36%  '
Set-TextValue 'D41' '0.1765'
Set-TextValue 'E41' '  -0.27%  '
Set-TextValue 'D42' '9.617'
Set-TextValue 'E42' '  -3.08%  '
Set-TextValue 'D43' '0.07088'
Set-TextValue 'E43' '  -3.58%  '
Set-TextValue 'D44' '11.49'
Set-TextValue 'E44' '  -1.64%  '
Set-TextValue 'D45' '0.5237'
Set-TextValue 'E45' '  -1.23%  '
Set-TextValue 'D46' '2.132'
Set-TextValue 'E46' '  -1.57%  '
Set-TextValue 'D47' '1.132'
Set-TextValue 'E47' '  -0.34%  '
Set-TextValue 'D48' '1.807'
Set-TextValue 'E48' '  -2.02%  '
Set-TextValue 'D49' '111.72'
Set-TextValue 'E49' '  -1.74%  '
Set-TextValue 'D50' '2.432'
Set-TextValue 'E50' '  +3.46%  '
Set-TextValue 'E51' '  +0.19%  '
